{"js": "// Find the \"Nomor:\" field placeholder (${no_sprin}) in the document body and\n// replace it with the fully-formatted Sprin letter-number template text,\n// keeping the same run formatting (Arial Narrow, bold-complex, sv-SE lang)\n// but split across several runs exactly as Word produced when the text was\n// typed/edited incrementally.\nconst searchResults = context.document.body.search(\"${no_sprin}\", { matchCase: true, matchWholeWord: false });\nsearchResults.load(\"items\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Could not find \"${no_sprin}\" placeholder in the document body.');\n}\n\nconst target = searchResults.items[0];\n\n// The new text is split into the same run boundaries seen in the authored\n// edit: \"Sprin.\" + \"Lidik${no_sprin}\" + \"/\" + \"XII\" + \"/HUK.6.6./202\" + \"2\"\n// which together read: Sprin.Lidik${no_sprin}/XII/HUK.6.6./2022\nconst pieces = [\"Sprin.\", \"Lidik${no_sprin}\", \"/\", \"XII\", \"/HUK.6.6./202\", \"2\"];\n\nconst runProps = '<w:rPr><w:rFonts w:ascii=\"Arial Narrow\" w:hAnsi=\"Arial Narrow\" w:cs=\"Arial\"/><w:bCs/><w:lang w:val=\"sv-SE\"/></w:rPr>';\nconst runsXml = pieces\n  .map((text) => `<w:r>${runProps}<w:t>${text}</w:t></w:r>`)\n  .join(\"\");\n\nconst ooxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>${runsXml}</w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\ntarget.insertOoxml(ooxml, Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Replace the \"Nomor: ${no_sprin}\" placeholder paragraph with the fully\n# expanded Sprin letter-number template text. The trailing placeholder run\n# is split into several runs (matching how Word records incremental typing)\n# that together read: Sprin.Lidik${no_sprin}/XII/HUK.6.6./2022\n# All other runs/formatting in the paragraph (\"No\" + \"mor\" + \": \") are left\n# exactly as they were.\n\n$d = $word.ActiveDocument\n\n$searchRange = $d.Content\n$find = $searchRange.Find\n$find.ClearFormatting()\n$find.Text = '${no_sprin}'\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$found = $find.Execute()\n\nif (-not $found) {\n    throw 'Could not find the \"${no_sprin}\" placeholder in the document.'\n}\n\n# Grow the found range out to the whole paragraph so we can replace it with\n# a fully reconstructed paragraph (InsertXML operates on whole paragraphs).\n$searchRange.Expand(4) | Out-Null  # 4 = wdParagraph\n\n$newParagraphXml = @'\n<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\" xmlns:w14=\"http://schemas.microsoft.com/office/word/2010/wordml\">\n        <w:body>\n          <w:p w14:paraId=\"1AD097B7\" w14:textId=\"45D78CF9\" w:rsidR=\"00476DAC\" w:rsidRPr=\"0040338A\" w:rsidRDefault=\"00476DAC\" w:rsidP=\"003415B9\">\n            <w:pPr>\n              <w:jc w:val=\"center\"/>\n              <w:rPr><w:rFonts w:ascii=\"Arial Narrow\" w:hAnsi=\"Arial Narrow\" w:cs=\"Arial\"/><w:bCs/><w:lang w:val=\"sv-SE\"/></w:rPr>\n            </w:pPr>\n            <w:r w:rsidRPr=\"0040338A\"><w:rPr><w:rFonts w:ascii=\"Arial Narrow\" w:hAnsi=\"Arial Narrow\" w:cs=\"Arial\"/><w:bCs/><w:lang w:val=\"sv-SE\"/></w:rPr><w:t>No</w:t></w:r>\n            <w:r w:rsidR=\"00E66C27\" w:rsidRPr=\"0040338A\"><w:rPr><w:rFonts w:ascii=\"Arial Narrow\" w:hAnsi=\"Arial Narrow\" w:cs=\"Arial\"/><w:bCs/><w:lang w:val=\"sv-SE\"/></w:rPr><w:t>mor</w:t></w:r>\n            <w:r w:rsidR=\"006167D6\" w:rsidRPr=\"0040338A\"><w:rPr><w:rFonts w:ascii=\"Arial Narrow\" w:hAnsi=\"Arial Narrow\" w:cs=\"Arial\"/><w:bCs/><w:lang w:val=\"sv-SE\"/></w:rPr><w:t xml:space=\"preserve\">: </w:t></w:r>\n            <w:r><w:rPr><w:rFonts w:ascii=\"Arial Narrow\" w:hAnsi=\"Arial Narrow\" w:cs=\"Arial\"/><w:bCs/><w:lang w:val=\"sv-SE\"/></w:rPr><w:t>Sprin.</w:t></w:r>\n            <w:r><w:rPr><w:rFonts w:ascii=\"Arial Narrow\" w:hAnsi=\"Arial Narrow\" w:cs=\"Arial\"/><w:bCs/><w:lang w:val=\"sv-SE\"/></w:rPr><w:t>Lidik${no_sprin}</w:t></w:r>\n            <w:r><w:rPr><w:rFonts w:ascii=\"Arial Narrow\" w:hAnsi=\"Arial Narrow\" w:cs=\"Arial\"/><w:bCs/><w:lang w:val=\"sv-SE\"/></w:rPr><w:t>/</w:t></w:r>\n            <w:r><w:rPr><w:rFonts w:ascii=\"Arial Narrow\" w:hAnsi=\"Arial Narrow\" w:cs=\"Arial\"/><w:bCs/><w:lang w:val=\"sv-SE\"/></w:rPr><w:t>XII</w:t></w:r>\n            <w:r><w:rPr><w:rFonts w:ascii=\"Arial Narrow\" w:hAnsi=\"Arial Narrow\" w:cs=\"Arial\"/><w:bCs/><w:lang w:val=\"sv-SE\"/></w:rPr><w:t>/HUK.6.6./202</w:t></w:r>\n            <w:r><w:rPr><w:rFonts w:ascii=\"Arial Narrow\" w:hAnsi=\"Arial Narrow\" w:cs=\"Arial\"/><w:bCs/><w:lang w:val=\"sv-SE\"/></w:rPr><w:t>2</w:t></w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>\n'@\n\n$searchRange.InsertXML($newParagraphXml)\n"}
